# upgrade Upload User function - lms
#
# Mirrors the sample row (row 2) of the "User List" sheet down into row 5,
# changes the sample e-mail value to "abcd", clears the stray leftover
# value that used to sit in K14, links the new row's e-mail cell the same
# way the existing sample rows are linked, stretches the Group dropdown
# down to the new row and leaves the saved cursor position on D7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User List")

# Clear the leftover scratch value first so its now-unused shared string
# slot gets reused (renamed) by the new value typed below, instead of a
# brand-new shared-string entry being appended.
$ws.Range("K14").ClearContents()

# Copy just the formatting of row 2 into row 5 (PasteSpecial "all" does not
# reproduce some of the per-cell styles reliably, so formats and values are
# applied separately).
$fmtCols = @("A","B","C","D","E","G","H")
foreach ($col in $fmtCols) {
    $ws.Range($col + "2").Copy()
    $ws.Range($col + "5").PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = 0

# Column F ("Enabled") already holds the literal text "true" under a Text
# number format; copy the whole cell so the value keeps its text type
# instead of being auto-coerced to a boolean.
$ws.Range("F2").Copy()
$ws.Range("F5").PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = 0

$ws.Range("A5").Value2 = "tien.duy2912@gmail.com"
$ws.Range("B5").Value2 = "tiennd_test0001"
$ws.Range("C5").Value2 = "Nguyen"
$ws.Range("D5").Value2 = "Tien"
$ws.Range("E5").Value2 = "abcd"
$ws.Range("G5").Value2 = "STUDENT"
$ws.Range("H5").Value2 = "STUDENT"

# Hyperlink the new row's e-mail cell the same way the other sample rows
# are linked.
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:tienduy2912@gmail.com", "", "", "tienduy2912@gmail.com") | Out-Null

# Move the saved selection/active cell.
$ws.Range("D7").Select()
